# Generate Report for Archive
#
# The localization status for these rows has moved on from the handoff
# step, so the "Ready for handoff" status text becomes "In Translation"
# everywhere it appears (Overview summary columns + each language
# sheet's Status column). Shrinking that text means the column that
# was sized to fit "Ready for handoff" is now too wide, so we resize
# the affected columns to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status text
$overview.Range("E2:F4").Value = "In Translation"
$overview.Columns("E:F").AutoFit() | Out-Null
$overview.Columns("E:F").ColumnWidth = 12.58

# zh-cn sheet: column C holds the status text
$zhcn.Range("C2:C4").Value = "In Translation"
$zhcn.Columns("C:C").AutoFit() | Out-Null
$zhcn.Columns("C:C").ColumnWidth = 12.58

# de-de sheet: column C holds the status text
$dede.Range("C2:C4").Value = "In Translation"
$dede.Columns("C:C").AutoFit() | Out-Null
$dede.Columns("C:C").ColumnWidth = 12.58
